# Refresh the cryptos price/volume list (Price column D and Volume(1h) column E),
# and for rows 46-48 also update the Coin name / Link (B/C) since the ranking order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value into column D while forcing a text format so that
# numeric-looking strings (e.g. "1.00", "0.0976") are preserved exactly as text
# instead of being auto-converted into floating point numbers by Excel. The
# cell style is reset back to Normal afterwards so the cell's appearance/style
# index stays the same as before the edit.
function Set-TextCell($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" '58.941.45'
$ws.Range("E2").Value = '  +0.12%  '
Set-TextCell $ws "D3" '2.568.98'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextCell $ws "D5" '568.49'
$ws.Range("E5").Value = '  +2.53%  '
Set-TextCell $ws "D6" '142.78'
$ws.Range("E6").Value = '  -0.48%  '
Set-TextCell $ws "D7" '1.00'
$ws.Range("E7").Value = '  +0.21%  '
Set-TextCell $ws "D8" '0.594'
$ws.Range("E8").Value = '  -0.20%  '
Set-TextCell $ws "D9" '2.570.23'
$ws.Range("E9").Value = '  -1.74%  '
Set-TextCell $ws "D10" '6.67'
$ws.Range("E10").Value = '  -2.19%  '
Set-TextCell $ws "D11" '0.103'
$ws.Range("E11").Value = '  +2.77%  '
$ws.Range("E12").Value = '  +11.16%  '
Set-TextCell $ws "D13" '0.343'
$ws.Range("E13").Value = '  +2.53%  '
Set-TextCell $ws "D14" '3.025.00'
$ws.Range("E14").Value = '  -1.21%  '
Set-TextCell $ws "D15" '58.994.08'
$ws.Range("E15").Value = '  +0.23%  '
Set-TextCell $ws "D16" '22.39'
$ws.Range("E16").Value = '  +7.36%  '
$ws.Range("E17").Value = '  +3.45%  '
Set-TextCell $ws "D18" '2.575.46'
$ws.Range("E18").Value = '  -1.40%  '
Set-TextCell $ws "D19" '4.52'
$ws.Range("E19").Value = '  +1.26%  '
Set-TextCell $ws "D20" '334.47'
$ws.Range("E20").Value = '  -0.94%  '
Set-TextCell $ws "D21" '10.20'
$ws.Range("E21").Value = '  +0.93%  '
Set-TextCell $ws "D22" '6.20'
$ws.Range("E22").Value = '  +0.55%  '
Set-TextCell $ws "D23" '1.00'
$ws.Range("E23").Value = '  +0.12%  '
Set-TextCell $ws "D24" '64.16'
$ws.Range("E24").Value = '  -3.45%  '
$ws.Range("E25").Value = '  +6.23%  '
Set-TextCell $ws "D26" '0.994'
$ws.Range("E26").Value = '  -0.22%  '
Set-TextCell $ws "D27" '0.160'
$ws.Range("E27").Value = '  +0.08%  '
Set-TextCell $ws "D28" '7.25'
$ws.Range("E28").Value = '  +0.90%  '
Set-TextCell $ws "D29" '0.0₃0776'
$ws.Range("E29").Value = '  +2.60%  '
Set-TextCell $ws "D30" '0.999'
$ws.Range("E30").Value = '  +0.14%  '
Set-TextCell $ws "D31" '1.68'
$ws.Range("E31").Value = '  +0.01%  '
Set-TextCell $ws "D32" '6.05'
$ws.Range("E32").Value = '  +0.93%  '
Set-TextCell $ws "D33" '158.33'
$ws.Range("E33").Value = '  +3.06%  '
Set-TextCell $ws "D34" '18.94'
$ws.Range("E34").Value = '  -0.25%  '
Set-TextCell $ws "D35" '4.03'
$ws.Range("E35").Value = '  +2.08%  '
Set-TextCell $ws "D36" '1.15'
$ws.Range("E36").Value = '  +1.54%  '
Set-TextCell $ws "D37" '0.869'
$ws.Range("E37").Value = '  -2.74%  '
Set-TextCell $ws "D38" '0.873'
$ws.Range("E38").Value = '  +0.33%  '
Set-TextCell $ws "D39" '37.06'
$ws.Range("E39").Value = '  +0.06%  '
Set-TextCell $ws "D40" '1.49'
$ws.Range("E40").Value = '  +2.31%  '
Set-TextCell $ws "D41" '3.66'
$ws.Range("E41").Value = '  +1.58%  '
Set-TextCell $ws "D42" '292.41'
$ws.Range("E42").Value = '  +3.64%  '
Set-TextCell $ws "D43" '0.999'
$ws.Range("E43").Value = '  +0.15%  '
Set-TextCell $ws "D44" '0.0976'
$ws.Range("E44").Value = '  +2.41%  '
Set-TextCell $ws "D45" '0.590'
$ws.Range("E45").Value = '  -1.70%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws "D46" '0.0535'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell $ws "D47" '10.63'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws "D48" '19.14'
$ws.Range("E48").Value = '  +1.72%  '
Set-TextCell $ws "D49" '125.04'
$ws.Range("E49").Value = '  +6.43%  '
Set-TextCell $ws "D50" '0.0232'
$ws.Range("E50").Value = '  +1.79%  '
Set-TextCell $ws "D51" '1.942.18'
$ws.Range("E51").Value = '  -0.12%  '
